$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells: "Throughput" -> "DLT" terminology, and the
# distance/size column header -> "File Size" (fairness-results rework).
$ws.Range("A1").Value = "File Size"
$ws.Range("B1").Value = "Avg. TCP DLT (Mbps)"
$ws.Range("C1").Value = "Avg. QUIC DLT (Mbps)"
